$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dataset for rows 2-31: Regiao (A), Ano (B), Valor (D), Posicao (E, Sergipe only)
$regioes = @(
  "Brasil",
  "Brasil",
  "Brasil",
  "Brasil",
  "Brasil",
  "Brasil",
  "Brasil",
  "Brasil",
  "Brasil",
  "Brasil",
  "Nordeste",
  "Nordeste",
  "Nordeste",
  "Nordeste",
  "Nordeste",
  "Nordeste",
  "Nordeste",
  "Nordeste",
  "Nordeste",
  "Nordeste",
  "Sergipe",
  "Sergipe",
  "Sergipe",
  "Sergipe",
  "Sergipe",
  "Sergipe",
  "Sergipe",
  "Sergipe",
  "Sergipe",
  "Sergipe"
)

$anos = @(
  "01/01/2015",
  "01/01/2016",
  "01/01/2017",
  "01/01/2018",
  "01/01/2019",
  "01/01/2020",
  "01/01/2021",
  "01/01/2022",
  "01/01/2023",
  "01/01/2024",
  "01/01/2015",
  "01/01/2016",
  "01/01/2017",
  "01/01/2018",
  "01/01/2019",
  "01/01/2020",
  "01/01/2021",
  "01/01/2022",
  "01/01/2023",
  "01/01/2024",
  "01/01/2015",
  "01/01/2016",
  "01/01/2017",
  "01/01/2018",
  "01/01/2019",
  "01/01/2020",
  "01/01/2021",
  "01/01/2022",
  "01/01/2023",
  "01/01/2024"
)

$valores = @(
  19.5367037911304,
  18.7691365547267,
  17.75241147846269,
  16.28100354585489,
  15.8148072839542,
  16.1193401823764,
  16.67367574713436,
  16.54723816555322,
  16.97629758730712,
  17.67162353660616,
  22.34782353990431,
  21.72884698728347,
  19.93727375877596,
  18.6593769057677,
  17.92702368561801,
  18.39897530445306,
  18.6127220282535,
  18.31537882657421,
  19.71373812688735,
  21.52258603851626,
  25.89068540205256,
  21.60748792710043,
  18.75518027117731,
  17.87103043104207,
  17.84678534778923,
  18.49865452232216,
  18.6760161540917,
  18.64892907205304,
  20.33309465318489,
  24.57359573685214
)

# Posicao relativamente as demais UF - only populated for the Sergipe rows (22-31)
$posicoes = @{
  22 = 9;
  23 = 14;
  24 = 16;
  25 = 15;
  26 = 15;
  27 = 15;
  28 = 16;
  29 = 14;
  30 = 11;
  31 = 7;
}

for ($i = 0; $i -lt $regioes.Length; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 1).Value = $regioes[$i]

  $cellB = $ws.Cells.Item($r, 2)
  $cellB.NumberFormat = "@"
  $cellB.Value = $anos[$i]
  $cellB.Style = "Normal"

  $ws.Cells.Item($r, 4).Value = $valores[$i]

  if ($posicoes.ContainsKey($r)) {
    $ws.Cells.Item($r, 5).Value = $posicoes[$r]
  }
}

# Drop the now-obsolete trailing Sergipe rows (old rows 32-34)
$ws.Rows.Item(32).Delete()
$ws.Rows.Item(32).Delete()
$ws.Rows.Item(32).Delete()

# Drop column F ("Faltam dados para todos os Estados") entirely
$ws.Columns.Item(6).Delete()
